# Fill in xG / goals data for the Benevento matches that were still
# missing it (rows 10-15: Benevento-Juventus, Parma-Benevento,
# Sassuolo-Benevento, Benevento-Lazio, Benevento-Genoa, Udinese-Benevento).
#
# The values look numeric but must be stored as TEXT (shared strings),
# matching the rest of the xG_home/xG_away/goals_home/goals_away columns.
# Setting NumberFormat="@" before assigning the Value keeps it text, and
# resetting the Style back to "Normal" afterwards drops the temporary
# number-format so no stray style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = 1
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D10" "0.439998"
Set-TextValue "E10" "1.64625"
Set-TextValue "F10" "1"
Set-TextValue "G10" "1"

Set-TextValue "D11" "0.501198"
Set-TextValue "E11" "0.340508"
Set-TextValue "F11" "0"
Set-TextValue "G11" "0"

Set-TextValue "D12" "1.08927"
Set-TextValue "E12" "2.09826"
Set-TextValue "F12" "1"
Set-TextValue "G12" "0"

Set-TextValue "D13" "1.43093"
Set-TextValue "E13" "0.416738"
Set-TextValue "F13" "1"
Set-TextValue "G13" "1"

Set-TextValue "D14" "1.55436"
Set-TextValue "E14" "0.245022"
Set-TextValue "F14" "2"
Set-TextValue "G14" "0"

Set-TextValue "D15" "2.4683"
Set-TextValue "E15" "0.338064"
Set-TextValue "F15" "0"
Set-TextValue "G15" "2"
